$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column P (2022) added alongside the existing year columns.
$ws.Range("P3").Value = 2022
$ws.Range("P4").Value = 15
$ws.Range("P5").Value = 2130.4

# Copy the formatting from the last existing column (O) onto the new one (P)
# so the new cells pick up the same style (borders, number format, font, etc.)
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial(-4122)

$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial(-4122)

[void]$ws.Range("P6").Select()
